# Generate Report for Handback
# -----------------------------------------------------------------------
# This script mirrors the "handback" report-generation step: the two
# localized files (zh-cn / de-de) have come back from translation and are
# now in sync with en-US, so:
#   1. every "Status" cell (and the mirrored Overview columns) flips from
#      "Ready for handoff" to "Handed back: in sync with en-US";
#   2. each language sheet gets its "Latest Target File" / "Latest
#      Handback File" / "Latest Handback DateTime" columns populated for
#      both rows, with a hyperlink added on the new "Latest Target File"
#      cell (mirroring the existing hyperlink on column A);
#   3. the widened columns are auto-fit so the longer strings are visible.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad898d9432628443c683ba84a26dc7721fd59101/e2e/"

$file1 = "06c2d3c9-32e1-4a2f-b211-d0c50b7828b0.md"
$file2 = "f31e9b55-f0c9-4503-a864-21f822572424.md"

# ---------------------------------------------------------------------
# Overview sheet: status columns for each language (E = zh-cn, F = de-de)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# Helper data: per-language sheet name, handback datetime and xlf names
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; Datetime = "2016-09-06 10:38:20";
       Xlf1 = "06c2d3c9-32e1-4a2f-b211-d0c50b7828b0.a6f440a71dd82e40d8398b64560380c9698886a8.zh-cn.xlf";
       Xlf2 = "f31e9b55-f0c9-4503-a864-21f822572424.5dae14079435160c6a080127ee1d5b59c4ee98c8.zh-cn.xlf" },
    @{ Sheet = "de-de"; Datetime = "2016-09-06 10:38:29";
       Xlf1 = "06c2d3c9-32e1-4a2f-b211-d0c50b7828b0.a6f440a71dd82e40d8398b64560380c9698886a8.de-de.xlf";
       Xlf2 = "f31e9b55-f0c9-4503-a864-21f822572424.5dae14079435160c6a080127ee1d5b59c4ee98c8.de-de.xlf" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) for both rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2 - 06c2d3c9-...md
    $ws.Range("I2").Value = $file1
    $ws.Range("J2").Value = $lang.Xlf1
    $ws.Range("K2").Value = $lang.Datetime

    # Row 3 - f31e9b55-...md
    $ws.Range("I3").Value = $file2
    $ws.Range("J3").Value = $lang.Xlf2
    $ws.Range("K3").Value = $lang.Datetime

    # New hyperlinks on the "Latest Target File" cells, mirroring column A
    $ws.Hyperlinks.Add($ws.Range("I2"), ($ghBase + $file1), "", "", $file1)
    $ws.Hyperlinks.Add($ws.Range("I3"), ($ghBase + $file2), "", "", $file2)

    # Give the new hyperlink cells the same look as the column-A hyperlinks
    # (underlined, cornflower-blue FF6495ED - matches the workbook's custom
    # "HyperLink" cell style already used on column A)
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = 15570276
    $ws.Range("I3").Font.Underline = $true
    $ws.Range("I3").Font.Color = 15570276

    # Widen the columns that now hold longer content
    $ws.Columns.Item(3).ColumnWidth = 29.15
    $ws.Columns.Item(9).ColumnWidth = 39.15
    $ws.Columns.Item(10).ColumnWidth = 39.15
}

Write-Host "Handback report generated"
